# Update cryptocurrency price (D) and volume-change (E) columns with
# the latest scraped values. Values are written as plain text (matching
# the source inline-string cells) by forcing a text number format for the
# assignment, then restoring the default "Normal" style so no extra
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '54.602.96'
Set-TextValue $ws.Range("E2") '  +1.36%  '

Set-TextValue $ws.Range("D3") '2.292.65'
Set-TextValue $ws.Range("E3") '  +2.04%  '

Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.12%  '

Set-TextValue $ws.Range("D5") '503.66'
Set-TextValue $ws.Range("E5") '  +2.42%  '

Set-TextValue $ws.Range("D6") '130.36'
Set-TextValue $ws.Range("E6") '  +2.89%  '

Set-TextValue $ws.Range("D7") '0.998'
Set-TextValue $ws.Range("E7") '  -0.12%  '

Set-TextValue $ws.Range("D8") '0.528'
Set-TextValue $ws.Range("E8") '  +0.78%  '

Set-TextValue $ws.Range("D9") '0.0963'
Set-TextValue $ws.Range("E9") '  +1.60%  '

Set-TextValue $ws.Range("E10") '  +0.89%  '

Set-TextValue $ws.Range("D11") '0.337'
Set-TextValue $ws.Range("E11") '  +4.48%  '

Set-TextValue $ws.Range("D12") '4.90'
Set-TextValue $ws.Range("E12") '  +6.08%  '

Set-TextValue $ws.Range("D13") '23.25'
Set-TextValue $ws.Range("E13") '  +7.42%  '

Set-TextValue $ws.Range("D14") '2.703.62'
Set-TextValue $ws.Range("E14") '  +1.81%  '

Set-TextValue $ws.Range("D15") '54.772.27'
Set-TextValue $ws.Range("E15") '  +1.69%  '

Set-TextValue $ws.Range("E16") '  +2.20%  '

Set-TextValue $ws.Range("D17") '2.306.76'
Set-TextValue $ws.Range("E17") '  +2.42%  '

Set-TextValue $ws.Range("D18") '10.35'
Set-TextValue $ws.Range("E18") '  +4.01%  '

Set-TextValue $ws.Range("D19") '4.17'
Set-TextValue $ws.Range("E19") '  +2.93%  '

Set-TextValue $ws.Range("D20") '306.98'
Set-TextValue $ws.Range("E20") '  +2.77%  '

Set-TextValue $ws.Range("D21") '6.35'
Set-TextValue $ws.Range("E21") '  -0.63%  '

Set-TextValue $ws.Range("D22") '1.00'
Set-TextValue $ws.Range("E22") '  +0.17%  '

Set-TextValue $ws.Range("D23") '60.80'
Set-TextValue $ws.Range("E23") '  -1.57%  '

Set-TextValue $ws.Range("D24") '0.991'
Set-TextValue $ws.Range("E24") '  -2.63%  '

Set-TextValue $ws.Range("D25") '0.150'
Set-TextValue $ws.Range("E25") '  +2.02%  '

Set-TextValue $ws.Range("D26") '7.42'
Set-TextValue $ws.Range("E26") '  +5.97%  '

Set-TextValue $ws.Range("D27") '173.20'
Set-TextValue $ws.Range("E27") '  +4.70%  '

Set-TextValue $ws.Range("D28") '0.0₃0724'
Set-TextValue $ws.Range("E28") '  +7.75%  '

Set-TextValue $ws.Range("D29") '6.05'
Set-TextValue $ws.Range("E29") '  +4.33%  '

Set-TextValue $ws.Range("D30") '1.62'
Set-TextValue $ws.Range("E30") '  +2.06%  '

Set-TextValue $ws.Range("E31") '  +4.97%  '

Set-TextValue $ws.Range("D33") '17.96'
Set-TextValue $ws.Range("E33") '  +2.28%  '

Set-TextValue $ws.Range("D34") '0.995'
Set-TextValue $ws.Range("E34") '  -0.13%  '

Set-TextValue $ws.Range("D35") '0.941'
Set-TextValue $ws.Range("E35") '  +6.58%  '

Set-TextValue $ws.Range("D36") '1.21'
Set-TextValue $ws.Range("E36") '  +3.43%  '

Set-TextValue $ws.Range("D37") '3.77'
Set-TextValue $ws.Range("E37") '  +3.81%  '

Set-TextValue $ws.Range("E38") '  +2.19%  '

Set-TextValue $ws.Range("D39") '1.42'
Set-TextValue $ws.Range("E39") '  +2.71%  '

Set-TextValue $ws.Range("D40") '3.41'
Set-TextValue $ws.Range("E40") '  +2.57%  '

Set-TextValue $ws.Range("D41") '4.88'
Set-TextValue $ws.Range("E41") '  -0.24%  '

Set-TextValue $ws.Range("D42") '125.62'
Set-TextValue $ws.Range("E42") '  +0.97%  '

Set-TextValue $ws.Range("D43") '249.55'
Set-TextValue $ws.Range("E43") '  +6.45%  '

Set-TextValue $ws.Range("D44") '0.0495'
Set-TextValue $ws.Range("E44") '  +3.31%  '

Set-TextValue $ws.Range("D45") '0.0899'
Set-TextValue $ws.Range("E45") '  +1.54%  '

Set-TextValue $ws.Range("D46") '0.550'
Set-TextValue $ws.Range("E46") '  +2.33%  '

Set-TextValue $ws.Range("D47") '0.376'
Set-TextValue $ws.Range("E47") '  +2.14%  '

Set-TextValue $ws.Range("D48") '0.0208'
Set-TextValue $ws.Range("E48") '  +3.76%  '

Set-TextValue $ws.Range("D49") '10.81'
Set-TextValue $ws.Range("E49") '  +0.71%  '

Set-TextValue $ws.Range("D50") '16.35'
Set-TextValue $ws.Range("E50") '  +2.86%  '

Set-TextValue $ws.Range("D51") '1.57'
Set-TextValue $ws.Range("E51") '  +7.29%  '
